# Update cryptos.xlsx price/volume/ranking data (GitHub Actions refresh)
# Source data refreshed from coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.846.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.21%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.084.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.67%  "
# Row 4
$ws.Range("E4").Value = "  +0.00%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.05%  "
# Row 8
$ws.Range("E8").Value = "  -0.04%  "
# Row 9
$ws.Range("E9").Value = "  -0.34%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0791"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.04%  "
# Row 11
$ws.Range("E11").Value = "  +3.03%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.392.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.69%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.39%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.20%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.769"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.19%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.91%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.083.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.63%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.756.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.00%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.34%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.54%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0834"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.88%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.91"
$ws.Range("D22").Style = "Normal"
# Row 23
$ws.Range("E23").Value = "  -0.05%  "
# Row 24
$ws.Range("E24").Value = "  -1.18%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.92%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.44%  "
# Row 27
$ws.Range("E27").Value = "  +7.51%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.97%  "
# Row 29
$ws.Range("E29").Value = "  +0.43%  "
# Row 30
$ws.Range("E30").Value = "  +2.35%  "
# Row 31
$ws.Range("E31").Value = "  +1.59%  "
# Row 32
$ws.Range("E32").Value = "  +3.62%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.29%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0630"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.31%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.27%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.00%  "
# Row 37
$ws.Range("E37").Value = "  +2.54%  "
# Row 38
$ws.Range("E38").Value = "  +0.06%  "
# Row 39
$ws.Range("E39").Value = "  -3.24%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0990"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.87%  "
# Row 41
$ws.Range("E41").Value = "  -0.06%  "
# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.88%  "
# Row 43
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.11%  "
# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.466.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.47%  "
# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0214"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.78%  "
# Row 46
$ws.Range("E46").Value = "  +0.86%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.12%  "
# Row 48
$ws.Range("E48").Value = "  +3.88%  "
# Row 49
$ws.Range("E49").Value = "  +2.55%  "
# Row 50
$ws.Range("E50").Value = "  +2.61%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.277.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.60%  "
